# Auto-generated: apply the cryptos-list price/volume refresh described in the diff.
# Column D ("Price") holds text-formatted values in the source data (e.g. "1.00",
# "43.100.78" with multiple separators) - force text format before writing so Excel
# does not silently coerce them into numbers and drop significant trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.100.78"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.08%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.288.69"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.72%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.50"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.95%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.531"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.21%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.506"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.60%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.20"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.45%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0818"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.70%  "

$ws.Range("E12").Value = "  +0.49%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.89"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.27%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.640.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.73%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.80"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.49%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.269.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.12%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.802"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.60%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.003.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0914"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.37%  "

$ws.Range("E21").Value = "  +0.98%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.85"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.66%  "

$ws.Range("E24").Value = "  +1.39%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.09%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.37%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.82"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.72%  "

$ws.Range("E29").Value = "  +1.22%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.30%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "165.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.25%  "

$ws.Range("E32").Value = "  +0.85%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.04%  "

$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.67%  "

$ws.Range("B35").Value = "Celestia"
$ws.Range("C35").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.04"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.25%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0738"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.42%  "

$ws.Range("E37").Value = "  +0.72%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.105"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.31%  "

$ws.Range("E39").Value = "  -0.08%  "

$ws.Range("E40").Value = "  +0.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.19"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.53%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.31"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.71%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.40"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.25%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0288"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.24%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.953.52"
$ws.Range("D45").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.45%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.30%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.86"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.97%  "

$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.54"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.83%  "

$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.509.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.93%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.76"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.94%  "
